$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix sheet name typo
$ws.Name = "ResgisterUsers"

# Insert a new first column (TCID) - shifts all existing data right by 1
[void]$ws.Columns("A").Insert()

# New TCID column - copy header formatting (bold/fill/border) from the
# neighboring header cell so the new column matches the rest of the header row
$ws.Range("B1").Copy() | Out-Null
[void]$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "TCID"
$ws.Range("A2").Value = 1

# Fix header text / rename headers (data-driven testing field names)
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "Repassword"
$ws.Range("E1").Value = "fname"
$ws.Range("F1").Value = "lname"
$ws.Range("G1").Value = "email"
$ws.Range("H1").Value = "address"
$ws.Range("I1").Value = "phone"

# Apply Text number format across the whole used range (matches the
# data-driven test's expectation that all fields are plain text)
$ws.Range("A1:K2").NumberFormat = "@"

# Move selection
[void]$ws.Range("H13").Select()

Write-Host "done"
